$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 2 (Ochieng Charles) to make room
# for "Jedidah Kemunto" as the new first data row.
$ws.Rows.Item(2).Insert()

# Insert two more rows after the (now shifted) Ochieng Charles row (row 3)
# for "Lenah Cheloti" and "Moses  Ngugi".
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

$dataRange = $ws.Range("A2:E6")

# Force text so numeric-looking strings (e.g. "1.00", "5.00%") are stored
# literally instead of being auto-converted to numbers/percentages, matching
# the source workbook's inlineStr cells.
$dataRange.NumberFormat = "@"

# Row 2: Jedidah Kemunto
$ws.Range("A2").Value = "Jedidah Kemunto"
$ws.Range("B2").Value = "0.00"
$ws.Range("C2").Value = "15.00"
$ws.Range("D2").Value = "-15.00"
$ws.Range("E2").Value = "0.00%"

# Row 3: Ochieng Charles (unchanged values, re-asserted since the row shifted)
$ws.Range("A3").Value = "Ochieng Charles"
$ws.Range("B3").Value = "1.00"
$ws.Range("C3").Value = "20.00"
$ws.Range("D3").Value = "-19.00"
$ws.Range("E3").Value = "5.00%"

# Row 4: Lenah Cheloti
$ws.Range("A4").Value = "Lenah Cheloti"
$ws.Range("B4").Value = "0.00"
$ws.Range("C4").Value = "12.00"
$ws.Range("D4").Value = "-12.00"
$ws.Range("E4").Value = "0.00%"

# Row 5: Moses  Ngugi
$ws.Range("A5").Value = "Moses  Ngugi"
$ws.Range("B5").Value = "0.00"
$ws.Range("C5").Value = "30.00"
$ws.Range("D5").Value = "-30.00"
$ws.Range("E5").Value = "0.00%"

# Row 6: KD Totals (updated totals)
$ws.Range("A6").Value = "KD Totals"
$ws.Range("B6").Value = "1.00"
$ws.Range("C6").Value = "77.00"
$ws.Range("D6").Value = "-76.00"
$ws.Range("E6").Value = "1.25%"

# Restore the default "Normal" style so the data rows carry no explicit
# cell style, matching the original workbook (row-insert otherwise
# inherits the bold/centred header style from row 1).
$dataRange.Style = "Normal"
